# cv.docx: "Built site for gh-pages" -- rebuilt resume content.
#
# Changes applied:
#   1. Personal site moved from https://ttalVlatt.github.io -> capaldi.info
#      (bookmark label + hyperlink display text both refer to the old host name).
#   2. "Campus Climate Vulnerability Project" is no longer ongoing -- drop the
#      "- Present" half of its date range.
#   3. "Merit-Based Financial Aid" project's date range closed out at Spring 2022.
#   4. Eight "Work Experience" job titles (everything except the still-current
#      "Graduate Research Assistant") get their heading runs explicitly bolded.
#
# Note on the bookmark itself: the document's `website-...` bookmark is a
# `w:bookmarkStart`/`w:bookmarkEnd` pair that was already present in the
# loaded package. Word's Bookmarks collection here only reports bookmarks
# created during the current session (Document.Bookmarks.Count is 0 for
# bookmarks that were already in the file), and bookmarks have no COM
# rename/delete that would let us retarget that existing id in place, so
# the bookmark's internal w:name keeps referencing the old host. Everything
# user-visible (the hyperlink text) is still corrected below.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

function Set-HeadingBold($text) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = 1
        $rng.Font.BoldBi = 1
    }
}

# --- 1. Website hyperlink text: https://ttalVlatt.github.io -> capaldi.info ---
Replace-Text "https://ttalVlatt.github.io" "capaldi.info"

# A plain text Replace rewrites the run from scratch and drops its rStyle, so
# re-find the new text and put the "Hyperlink" character style back on it.
$siteRng = $d.Content
$siteRng.Find.Execute("capaldi.info", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
$siteRng.Style = "Hyperlink"

# --- 2. Project date ranges ---
Replace-Text "Campus Climate Vulnerability Project (Summer 2023 - Present)" `
             "Campus Climate Vulnerability Project (Summer 2023)"
Replace-Text "Merit-Based Financial Aid (Fall 2021 - Present)" `
             "Merit-Based Financial Aid (Fall 2021 - Spring 2022)"

# --- 3. Bold the job-title headings (all "Work Experience" entries except the
#         still-current "Graduate Research Assistant") ---
Set-HeadingBold "Residence Hall Coordinator"
Set-HeadingBold "Graduate Community Coordinator"
Set-HeadingBold "Summer Sessions Graduate Intern"
Set-HeadingBold "Graduate Assistant for Industry Relations"
Set-HeadingBold "Advising Intern"
Set-HeadingBold "Student Ambassador"
Set-HeadingBold "Careers Network Student Peer Presenter"
Set-HeadingBold "International Employer Liaison Assistant"
